$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.858.84"
$ws.Range("E2").Value = "  +4.83%  "
$ws.Range("D3").Value = "2.974.84"
$ws.Range("E3").Value = "  +2.39%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'580.36"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.93%  "
$ws.Range("D6").Value = "'153.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +6.33%  "
$ws.Range("E7").Value = "  +0.15%  "
$ws.Range("D8").Value = "2.968.11"
$ws.Range("E8").Value = "  +2.18%  "
$ws.Range("E9").Value = "  +1.02%  "
$ws.Range("D10").Value = "'6.98"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.15%  "
$ws.Range("E11").Value = "  +2.03%  "
$ws.Range("E12").Value = "  +2.57%  "
$ws.Range("E13").Value = "  +2.88%  "
$ws.Range("D14").Value = "'34.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +6.48%  "
$ws.Range("E15").Value = "  +0.68%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "64.759.24"
$ws.Range("E16").Value = "  +4.73%  "
$ws.Range("B17").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C17").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D17").Value = "3.472.62"
$ws.Range("E17").Value = "  +2.71%  "
$ws.Range("D18").Value = "'6.91"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.85%  "
$ws.Range("D19").Value = "2.973.63"
$ws.Range("E19").Value = "  +2.73%  "
$ws.Range("D20").Value = "'447.39"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.30%  "
$ws.Range("D21").Value = "'13.58"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.24%  "
$ws.Range("D22").Value = "'0.676"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.71%  "
$ws.Range("D23").Value = "'7.24"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.58%  "
$ws.Range("D24").Value = "'80.87"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.99%  "
$ws.Range("D25").Value = "'11.02"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +8.24%  "
$ws.Range("D26").Value = "'12.22"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.42%  "
$ws.Range("E27").Value = "  +7.37%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("E29").Value = "  +9.53%  "
$ws.Range("D30").Value = "'0.0000108"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.97%  "
$ws.Range("D31").Value = "'2.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.93%  "
$ws.Range("E32").Value = "  +1.97%  "
$ws.Range("E33").Value = "  +2.16%  "
$ws.Range("D34").Value = "'26.51"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.12%  "
$ws.Range("D35").Value = "'0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("D36").Value = "'0.979"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.16%  "
$ws.Range("D37").Value = "'5.63"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.96%  "
$ws.Range("D38").Value = "'2.12"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.19%  "
$ws.Range("E39").Value = "  +6.09%  "
$ws.Range("D40").Value = "'48.81"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.52%  "
$ws.Range("D41").Value = "'44.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +10.60%  "
$ws.Range("D42").Value = "'0.119"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.85%  "
$ws.Range("D43").Value = "'0.295"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +8.99%  "
$ws.Range("D44").Value = "'8.34"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.40%  "
$ws.Range("D45").Value = "'388.32"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +14.25%  "
$ws.Range("D46").Value = "2.777.15"
$ws.Range("E46").Value = "  +3.11%  "
$ws.Range("E47").Value = "  +4.90%  "
$ws.Range("D48").Value = "'135.26"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.97%  "
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").Value = "'0.000226"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +14.76%  "
$ws.Range("E51").Value = "  +1.87%  "
